## Added Week 15 simulations
## Insert a new "N.Mullens" sheet before the existing "B.Mayfield" sheet,
## mirroring the layout/content of the other QB sheets (headers + H/R rows
## of zeroed-out stats).

$wb = $excel.ActiveWorkbook

# Existing first sheet (B.Mayfield) - the new sheet goes in front of it so
# the final tab order becomes: N.Mullens, B.Mayfield, C.Keenum
$firstSheet = $wb.Worksheets.Item(1)

$ws = $wb.Worksheets.Add($firstSheet)
$ws.Name = "N.Mullens"

# Column headers (row 1, B:G) - reuses the same shared strings as the
# other player sheets.
$ws.Range("B1").Value = "Short Att"
$ws.Range("C1").Value = "Short Comp"
$ws.Range("D1").Value = "Deep Att"
$ws.Range("E1").Value = "Deep Comp"
$ws.Range("F1").Value = "Short Int"
$ws.Range("G1").Value = "Deep Int"

# Row labels (Home / Road)
$ws.Range("A2").Value = "H"
$ws.Range("A3").Value = "R"

# Data cells start at zero (simulation not yet populated)
$ws.Range("B2:G3").Value = 0

# Match the bold / centered / bordered header style used on the other
# sheets for the header row and the row labels.
foreach ($rng in @($ws.Range("B1:G1"), $ws.Range("A2:A3"))) {
    $rng.Font.Bold = $true
    $rng.Borders.LineStyle = 1
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
}
